$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for existing rows 2-28
# from 45192 (2023-09-23) to 45202 (2023-10-03).
for ($r = 2; $r -le 28; $r++) {
    $ws.Range("C$r").Value = 45202
}

# Row 28 becomes an explicit custom-height row (height 15) once a new
# row is appended after it.
$ws.Rows.Item(28).RowHeight = 15

# Add the new data row 29.
$ws.Range("A29").Value = "A 46947-2023"
$ws.Range("B29").Value = 45201
$ws.Range("C29").Value = 45202
$ws.Range("D29").Value = "SKÅNE LÄN"
$ws.Range("E29").Value = "BÅSTAD"
$ws.Range("G29").Value = 1.5
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0

# R29 keeps the wrap-text style used by the rest of column R, with no content.
$ws.Range("R29").WrapText = $true

# Apply the date number format (style used by B/C columns) to the new date cells.
$ws.Range("B29").NumberFormat = $ws.Range("B28").NumberFormat()
$ws.Range("C29").NumberFormat = $ws.Range("C28").NumberFormat()
